$d = $word.ActiveDocument
$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Merge the two runs in the "Which clients do not have backups (onprem/
#    cloud)" bullet into a single run.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Which clients do not have backups*") {
        $target = $p
        break
    }
}
$xml = '<w:p ' + $wordmlNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Which clients do not have backups (onprem/cloud)</w:t></w:r></w:p>'
$target.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Drop the _GoBack bookmark that trails the "Employees need..." bullet.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Employees need to be able*") {
        $target = $p
        break
    }
}
$xml = '<w:p ' + $wordmlNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Employees need to be able to access the system and view/modify according to their position</w:t></w:r></w:p>'
$target.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) Append the three new normalization bullet points after the picture
#    paragraph (the document's last paragraph), ahead of the sectPr.
# ---------------------------------------------------------------------------
$picturePara = $d.Paragraphs.Last

$picturePara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($picturePara.Index + 1)
$xml1 = '<w:p ' + $wordmlNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>The data is already normalized.</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

$p1 = $d.Paragraphs.Item($picturePara.Index + 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($picturePara.Index + 2)
$xml2 = '<w:p ' + $wordmlNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">It is first normal form because there are no repeating columns, or values within the data. </w:t></w:r></w:p>'
$p2.Range.InsertXML($xml2)

$p2 = $d.Paragraphs.Item($picturePara.Index + 2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($picturePara.Index + 3)
$xml3 = '<w:p ' + $wordmlNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>It is second</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> and third</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> normal form because all non-key values depend completely</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> and only</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> on the primary key.</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p3.Range.InsertXML($xml3)

Write-Output "done"
